$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three data values on row 1
$ws.Range("A1").Value = 161.63079833413553
$ws.Range("B1").Value = 6.7983330273911902
$ws.Range("C1").Value = 0.79373776908023486

# Widen column C by one character unit (11.7109375 -> 12.7109375)
$ws.Columns.Item(3).ColumnWidth = 11.8
